$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.298.21'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '3.496.06'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'588.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").Value = "'133.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").Value = "'7.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.28%  '
$ws.Range("E10").Value = '  -0.40%  '
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("D12").Value = '4.094.80'
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("E14").Value = '  -1.12%  '
$ws.Range("D15").Value = '3.501.73'
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("D16").Value = '64.250.43'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = "'24.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.43%  '
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("E20").Value = '  -1.85%  '
$ws.Range("D21").Value = "'386.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.25%  '
$ws.Range("E22").Value = '  +1.87%  '
$ws.Range("D23").Value = '3.635.54'
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = "'74.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("E27").Value = '  +0.95%  '
$ws.Range("D28").Value = "'0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("E30").Value = '  +0.28%  '
$ws.Range("E31").Value = '  +0.84%  '
$ws.Range("D32").Value = "'8.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.00%  '
$ws.Range("E33").Value = '  +3.38%  '
$ws.Range("D34").Value = '3.524.08'
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("D36").Value = "'23.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.35%  '
$ws.Range("D37").Value = "'5.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.42%  '
$ws.Range("D38").Value = "'6.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("D39").Value = "'1.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.51%  '
$ws.Range("D40").Value = "'164.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.71%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").Value = "'0.806"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("E44").Value = '  -0.71%  '
$ws.Range("E45").Value = '  +0.82%  '
$ws.Range("D46").Value = "'24.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.44%  '
$ws.Range("E47").Value = '  -0.82%  '
$ws.Range("D48").Value = '2.425.65'
$ws.Range("E48").Value = '  -2.36%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = "'6.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").Value = "'0.919"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.23%  '
$ws.Range("E51").Value = '  -1.58%  '
